# Update the workbook for the 2022-08-08 data refresh (commit: "Add data for 2022-08-16")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-08"

# Update the header label in I1 (shared string "2022 (through 08-07)" -> "2022 (through 08-08)")
$ws.Range("I1").Value = "2022 (through 08-08)"

# Update the September value and the Total value in column I
$ws.Range("I9").Value = 42
$ws.Range("I14").Value = 1012
